# Update "Pais" sheet: refresh COVID-19 country stats and reorder Barein row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "last updated" timestamp text (cell A1)
$ws.Range("A1").Value = "Datos actualizados a 18 de Abril de 2020 a las 12:52"

# 2. Iran (row 11)
$ws.Cells.Item(11,2).Value = 80868   # Casos totales
$ws.Cells.Item(11,3).Value = 1374    # Nuevos casos
$ws.Cells.Item(11,4).Value = 55987   # Casos activos
$ws.Cells.Item(11,5).Value = 19850   # Recuperados
$ws.Cells.Item(11,6).Value = 3513    # Casos criticos
$ws.Cells.Item(11,7).Value = 73      # Muertes hoy
$ws.Cells.Item(11,8).Value = 5031    # Muertes

# 3. Suiza (row 18)
$ws.Cells.Item(18,5).Value = 9334    # Recuperados
$ws.Cells.Item(18,7).Value = 17      # Muertes hoy
$ws.Cells.Item(18,8).Value = 1344    # Muertes

# 4. Austria (row 20)
$ws.Cells.Item(20,2).Value = 14671   # Casos totales
$ws.Cells.Item(20,3).Value = 76      # Nuevos casos
$ws.Cells.Item(20,5).Value = 4014    # Recuperados
$ws.Cells.Item(20,7).Value = 12      # Muertes hoy
$ws.Cells.Item(20,8).Value = 443     # Muertes

# 5. Finlandia (row 50)
$ws.Cells.Item(50,2).Value = 3681    # Casos totales
$ws.Cells.Item(50,3).Value = 192     # Nuevos casos
$ws.Cells.Item(50,5).Value = 1891    # Recuperados
$ws.Cells.Item(50,7).Value = 8       # Muertes hoy
$ws.Cells.Item(50,8).Value = 90      # Muertes

# 6. Reorder Barein ahead of Islandia/Kuwait (rows 64-66) with refreshed data
#    Row 64 becomes Barein (new figures), Islandia and Kuwait shift down one row
#    keeping their previous figures unchanged.
$ws.Cells.Item(64,1).Value = "Barein"
$ws.Cells.Item(64,2).Value = 1767
$ws.Cells.Item(64,3).Value = 27
$ws.Cells.Item(64,4).Value = 741
$ws.Cells.Item(64,5).Value = 1019
$ws.Cells.Item(64,6).Value = 3
$ws.Cells.Item(64,7).Value = 0
$ws.Cells.Item(64,8).Value = 7

$ws.Cells.Item(65,1).Value = "Islandia"
$ws.Cells.Item(65,2).Value = 1754
$ws.Cells.Item(65,3).Value = 0
$ws.Cells.Item(65,4).Value = 1224
$ws.Cells.Item(65,5).Value = 521
$ws.Cells.Item(65,6).Value = 6
$ws.Cells.Item(65,7).Value = 0
$ws.Cells.Item(65,8).Value = 9

$ws.Cells.Item(66,1).Value = "Kuwait"
$ws.Cells.Item(66,2).Value = 1751
$ws.Cells.Item(66,3).Value = 93
$ws.Cells.Item(66,4).Value = 280
$ws.Cells.Item(66,5).Value = 1465
$ws.Cells.Item(66,6).Value = 34
$ws.Cells.Item(66,7).Value = 1
$ws.Cells.Item(66,8).Value = 6

# 7. Gibraltar (row 133)
$ws.Cells.Item(133,4).Value = 111    # Casos activos
$ws.Cells.Item(133,5).Value = 21     # Recuperados

# 8. Trinidad yTobago (row 137)
$ws.Cells.Item(137,4).Value = 21     # Casos activos
$ws.Cells.Item(137,5).Value = 85     # Recuperados
